# Rougelike RPG.docx -- apply the target edit via Word COM-interop calls.
#
# Target shape of the document after the edit (3 paragraphs total):
#   1) "You are the son of a legendary dungeoneer who's gone missing for a
#       year. You are determined to enter the dungeon yourself and find out
#       what has happened to him. " + a trailing " " run.
#   2) (unchanged) empty paragraph holding the _GoBack bookmark.
#   3) "Our story pillars: Character and Lore. Each <npc> will have some
#       lines to say, and a cool design. " (keeping the proofing wrapper
#       around "npc" that used to wrap "ai").
# Every other paragraph in the original doc is removed.

$d = $word.ActiveDocument

# --- Step 1: drop every paragraph after the "ai" paragraph (originally
#     paragraph 5) through the end of the document. ---
$tailStart = $d.Paragraphs.Item(6).Range.Start
$tailEnd = $d.Paragraphs.Item($d.Paragraphs.Count).Range.End
$d.Range($tailStart, $tailEnd).Delete()

# --- Step 2: merge paragraph 3 ("Focus on making...") together with the
#     blank paragraph 4 and paragraph 5 ("The skill-based gameplay...") by
#     removing the two intervening paragraph marks one at a time (Word
#     joins the runs of the paragraph being absorbed into the following
#     paragraph when its mark is deleted). ---
$mark1 = $d.Paragraphs.Item(3).Range.End - 1
$d.Range($mark1, $mark1 + 1).Delete()
$mark2 = $d.Paragraphs.Item(3).Range.End - 1
$d.Range($mark2, $mark2 + 1).Delete()

# --- Step 3: rewrite paragraph 1. ---
$p1 = $d.Paragraphs.Item(1).Range
$p1.Find.Execute(
    "Overall Goal: Create an arena brawler ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "You are the son of a legendary dungeoneer who" + [char]8217 + "s gone missing for a year. You are determined to enter the dungeon yourself and find out what has happened to him. ",
    2) | Out-Null

$p1b = $d.Paragraphs.Item(1).Range
$p1b.Find.Execute(
    "rpg, focusing on immersing the player in the progression, and the skill-based gameplay. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " ",
    2) | Out-Null

# --- Step 4: rewrite the merged paragraph 3. ---
# 4a: first run's text becomes the new story-pillars intro.
$p3 = $d.Paragraphs.Item(3).Range
$p3.Find.Execute(
    "Focus on making an exciting and involving progression: Items, Skills, Progression Trees. I am hoping that the progression tree will be dynamic, where the player must work towards a skill to ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Our story pillars: Character and Lore. Each ",
    2) | Out-Null

# 4b: drop everything between "Each " and the first occurrence of "ai".
$p3b = $d.Paragraphs.Item(3).Range
$midRange1 = $d.Range($p3b.Start, $p3b.End)
$midRange1.Find.Execute(
    "access certain sections of the skills tree. I am thinking of not putting any numbers on the equipment, and the only way to test is to do combat, or to go in the dummy arena. I also think the items should involve a good deal of combinatorics, where certain combinations will work better than others. We have to be very careful about balancing however. The skill-based gameplay will hinge on good ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0) | Out-Null
$midRange1.Delete()

# 4c: rename the first (still proof-wrapped) "ai" occurrence to "npc".
$p3c = $d.Paragraphs.Item(3).Range
$p3c.Find.Execute(
    "ai", $true, $true, $false, $false, $false, $true, 1, $false,
    "npc", 1) | Out-Null

# 4d: collapse the remainder (second "ai" occurrence + its proofing marks)
#     down to the final sentence.
$p3d = $d.Paragraphs.Item(3).Range
$midRange2 = $d.Range($p3d.Start, $p3d.End)
$midRange2.Find.Execute(
    ", and fun/effective/well-balanced skills. The ai must be able to use skills as well. Some skills should be obvious when to use, and some should be less obvious, but more rewarding when it does work against an enemy. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    " will have some lines to say, and a cool design. ",
    2) | Out-Null
